$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '38.418.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.086.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.75'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.398.64'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.38'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.787'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.088.57'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '38.332.85'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.44'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.48'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.78'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.136'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.13%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.03%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.57%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.81'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.41%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.34'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.537.93'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.88'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0220'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0931'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.69'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +8.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.283.29'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.02%  '
